$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Title (appears twice - main heading and bold recap near end)
Replace-Text "Play Astro Legends for Free: Review of Lyra and Erion Online Slot" "Play Astro Legends: Lyra and Erion Free - Slot Game Review"

# "What we like" bullet list
Replace-Text "Hexagonal, 9-reel gameplay with diverse winning combinations" "Hexagonal gameplay with 9 reels"
Replace-Text "Sonic Respin bonus for more wins and free spins" "Respin function for improved winnings"
Replace-Text "Valuable wild multipliers up to 5x the bet" "Wide betting range from 0.10 cents to 100.00 euros"
Replace-Text "Lyra Spirit Bonus with 7 levels for bigger rewards" "Sonic Respin Bonus and Lyra Spirit Bonus for extra excitement"

# "What we don't like" bullet list
Replace-Text "Gameplay might be overwhelming for some players due to its complexity" "Limited number of symbols in the winning group"
Replace-Text "Players may experience an occasional lag or delay while playing" "Chance of losing a life in the Lyra Spirit Bonus round"

# Meta description (italic) near the end
Replace-Text "Read our review of Astro Legends: Lyra and Erion online slot. Play for free at SlotJava platform. Find out about the Sonic Respin, Lyra Spirit Bonus, and more." "Discover the thrilling gameplay of Astro Legends: Lyra and Erion in this free slot game review."
